# Applies the update described by the commit:
#  - Column C ("Förändrad") date bumped from 45208 to 45212 for all data rows (2-108)
#  - For the first two data rows (2 and 3), the hyperlink formulas in columns S-Y
#    get updated target filenames (adding descriptive suffixes to the linked files)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C (Förändrad) for all data rows 2..108 ---
for ($r = 2; $r -le 108; $r++) {
    $ws.Cells.Item($r, 3).Value = 45212
}

# --- Update hyperlink formulas for row 2 (A 30779-2023 / Logging_0883) ---
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/artfynd/A 30779-2023 artfynd.xlsx", "A 30779-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/kartor/A 30779-2023 karta.png", "A 30779-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/knärot/A 30779-2023 karta knärot.png", "A 30779-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/klagomål/A 30779-2023 fsc-klagomål.docx", "A 30779-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/klagomålsmail/A 30779-2023 fsc-klagomål mail.docx", "A 30779-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/tillsyn/A 30779-2023 tillsynsbegäran.docx", "A 30779-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/ti,llsynsmail/A 30779-2023 tillsynsbegäran mail.docx", "A 30779-2023")'

# --- Update hyperlink formulas for row 3 (A 32298-2023 / Logging_0861) ---
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/artfynd/A 32298-2023 artfynd.xlsx", "A 32298-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/kartor/A 32298-2023 karta.png", "A 32298-2023")'
$ws.Range("U3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/knärot/A 32298-2023 karta knärot.png", "A 32298-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/klagomål/A 32298-2023 fsc-klagomål.docx", "A 32298-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/klagomålsmail/A 32298-2023 fsc-klagomål mail.docx", "A 32298-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/tillsyn/A 32298-2023 tillsynsbegäran.docx", "A 32298-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/ti,llsynsmail/A 32298-2023 tillsynsbegäran mail.docx", "A 32298-2023")'
